$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - human readable column headers
$ws.Range("B1").Value = "Grandes grupos, código"
$ws.Range("C1").Value = "Provincia codigo"
$ws.Range("D1").Value = "Provincia:"
$ws.Range("E1").Value = "Comarca nombre"
$ws.Range("F1").Value = "Total"
$ws.Range("G1").Value = "Sexo, código"
$ws.Range("H1").Value = "Españoles"
$ws.Range("I1").Value = "Extranjeros"
$ws.Range("J1").Value = "Sexo"
$ws.Range("K1").Value = "Grandes grupos"
$ws.Range("L1").Value = "Comarca codigo"

# Row 2 - measure / dimension identifiers
$ws.Range("A2").Value = "iaest-measure:"
$ws.Range("B2").Value = "null"
$ws.Range("C2").Value = "null"
$ws.Range("D2").Value = "iaest-measure:provincia"
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "iaest-measure:total"
$ws.Range("G2").Value = "null"
$ws.Range("H2").Value = "iaest-measure:espanoles"
$ws.Range("I2").Value = "iaest-measure:extranjeros"
$ws.Range("J2").Value = "iaest-measure:sexo"
$ws.Range("K2").Value = "iaest-measure:grandes-grupos"
$ws.Range("L2").Value = "null"

# Row 3 - medida/dim markers
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "null"
$ws.Range("C3").Value = "null"
$ws.Range("D3").Value = "medida"
$ws.Range("E3").Value = "dim"
$ws.Range("F3").Value = "medida"
$ws.Range("G3").Value = "null"
$ws.Range("H3").Value = "medida"
$ws.Range("I3").Value = "medida"
$ws.Range("J3").Value = "medida"
$ws.Range("K3").Value = "medida"
$ws.Range("L3").Value = "null"

# Row 4 - data types / URI
$ws.Range("A4").Value = "xsd:double"
$ws.Range("B4").Value = "null"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "xsd:string"
$ws.Range("E4").Value = "URI-comarca"
$ws.Range("F4").Value = "xsd:double"
$ws.Range("G4").Value = "null"
$ws.Range("H4").Value = "xsd:double"
$ws.Range("I4").Value = "xsd:double"
$ws.Range("J4").Value = "xsd:string"
$ws.Range("K4").Value = "xsd:string"
$ws.Range("L4").Value = "null"
